$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text / shared-string updates -----------------------------------------

# Police Commissioner name (plain string cell)
$ws.Range("M6").Value = "Thomas G. Donlon"

# Volume/Number line: "Volume 31   Number  38" -> "...39"
$numChars = $ws.Range("A8").Characters(21, 2)
$numChars.Text = "39"

# Report-covering-week line: dates shift forward one week
$d1 = $ws.Range("C9").Characters(27, 9)
$d1.Text = "9/23/2024"
$d2 = $ws.Range("C9").Characters(47, 9)
$d2.Text = "9/29/2024"

# --- Column widths ----------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 7.433768
$ws.Columns.Item(8).ColumnWidth = 6.168446

# --- Data rows ---------------------------------------------------------------

$ws.Range("N15").Value = -52.380952380952

$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 7
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 80
$ws.Range("J16").Value = 104
$ws.Range("K16").Value = -23.076923076923
$ws.Range("L16").Value = -28.571428571428
$ws.Range("M16").Value = 2.564102564102
$ws.Range("N16").Value = -82.832618025751

$ws.Range("C17").Value = 5
$ws.Range("E17").Value = 150
$ws.Range("F17").Value = 10
$ws.Range("H17").Value = -9.090909090909
$ws.Range("I17").Value = 95
$ws.Range("J17").Value = 123
$ws.Range("K17").Value = -22.764227642276
$ws.Range("L17").Value = -5
$ws.Range("M17").Value = -2.061855670103
$ws.Range("N17").Value = -48.369565217391

$ws.Range("C18").Value = 0
$ws.Range("F18").Value = 4
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = -20
$ws.Range("L18").Value = -45.864661654135
$ws.Range("M18").Value = -11.111111111111
$ws.Range("N18").Value = -84.210526315789

$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 28.571428571428
$ws.Range("F19").Value = 54
$ws.Range("G19").Value = 41
$ws.Range("H19").Value = 31.707317073170
$ws.Range("I19").Value = 516
$ws.Range("J19").Value = 570
$ws.Range("K19").Value = -9.473684210526
$ws.Range("L19").Value = -1.901140684410
$ws.Range("M19").Value = 10.256410256410
$ws.Range("N19").Value = -16.233766233766

$ws.Range("C20").Value = 0
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = "***.*"
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = -20
$ws.Range("M20").Value = 14.285714285714
$ws.Range("N20").Value = -89.473684210526

$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 10
$ws.Range("E21").Value = 60
$ws.Range("F21").Value = 83
$ws.Range("G21").Value = 74
$ws.Range("H21").Value = 12.162162162162
$ws.Range("I21").Value = 813
$ws.Range("J21").Value = 959
$ws.Range("K21").Value = -15.224191866527
$ws.Range("L21").Value = -13.877118644067
$ws.Range("M21").Value = 6.135770234986
$ws.Range("N21").Value = -61.795112781954

$ws.Range("C22").Value = 2
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 23
$ws.Range("K22").Value = -11.538461538461
$ws.Range("L22").Value = 4.545454545454
$ws.Range("M22").Value = 155.555555555556

$ws.Range("D23").Value = 2
$ws.Range("E23").Value = -50
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = -60
$ws.Range("I23").Value = 32
$ws.Range("J23").Value = 38
$ws.Range("K23").Value = -15.789473684210
$ws.Range("L23").Value = -25.581395348837
$ws.Range("M23").Value = -3.030303030303

$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = -31.578947368421
$ws.Range("G24").Value = 56
$ws.Range("H24").Value = -1.785714285714
$ws.Range("I24").Value = 511
$ws.Range("J24").Value = 617
$ws.Range("K24").Value = -17.179902755267
$ws.Range("L24").Value = -8.586762075134
$ws.Range("M24").Value = -22.103658536585

$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 200
$ws.Range("F25").Value = 26
$ws.Range("G25").Value = 20
$ws.Range("H25").Value = 30
$ws.Range("I25").Value = 226
$ws.Range("J25").Value = 355
$ws.Range("K25").Value = -36.338028169014
$ws.Range("L25").Value = -28.930817610062

$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 100
$ws.Range("F26").Value = 31
$ws.Range("G26").Value = 26
$ws.Range("H26").Value = 19.230769230769
$ws.Range("I26").Value = 243
$ws.Range("J26").Value = 254
$ws.Range("K26").Value = -4.330708661417
$ws.Range("L26").Value = -15.625
$ws.Range("M26").Value = -5.813953488372

$ws.Range("C27").Value = 0

$ws.Range("C28").Value = 0
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 6
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 20
$ws.Range("J28").Value = 47
$ws.Range("K28").Value = -2.127659574468
$ws.Range("L28").Value = 4.545454545454

$ws.Range("D31").Value = 1
$ws.Range("E31").Value = -100
$ws.Range("F31").Value = 2
$ws.Range("H31").Value = -33.333333333333
$ws.Range("I31").Value = 7
$ws.Range("J31").Value = 17
$ws.Range("K31").Value = -58.823529411764
$ws.Range("L31").Value = -56.25
